# Auto-generated edit script to update F-column ('想去人数' / want-to-go counts)
# across multiple worksheets, per the commit's output regeneration diff.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Cells.Item(2, 6).Value = 1562
$ws.Cells.Item(4, 6).Value = 2114
$ws.Cells.Item(5, 6).Value = 8722
$ws.Cells.Item(6, 6).Value = 234
$ws.Cells.Item(7, 6).Value = 96
$ws.Cells.Item(8, 6).Value = 1240
$ws.Cells.Item(10, 6).Value = 239
$ws.Cells.Item(11, 6).Value = 591
$ws.Cells.Item(13, 6).Value = 103
$ws.Cells.Item(14, 6).Value = 278
$ws.Cells.Item(16, 6).Value = 46
$ws.Cells.Item(17, 6).Value = 1402
$ws.Cells.Item(18, 6).Value = 1305
$ws.Cells.Item(19, 6).Value = 564
$ws.Cells.Item(21, 6).Value = 1324
$ws.Cells.Item(22, 6).Value = 66
$ws.Cells.Item(23, 6).Value = 205
$ws.Cells.Item(25, 6).Value = 64
$ws.Cells.Item(26, 6).Value = 57
$ws.Cells.Item(27, 6).Value = 273
$ws.Cells.Item(28, 6).Value = 1057
$ws.Cells.Item(29, 6).Value = 5
$ws.Cells.Item(31, 6).Value = 199
$ws.Cells.Item(32, 6).Value = 174
$ws.Cells.Item(35, 6).Value = 596
$ws.Cells.Item(37, 6).Value = 117
$ws.Cells.Item(39, 6).Value = 142
$ws.Cells.Item(40, 6).Value = 463
$ws.Cells.Item(42, 6).Value = 664
$ws.Cells.Item(43, 6).Value = 195
$ws.Cells.Item(45, 6).Value = 38

$ws = $wb.Worksheets.Item("演出")
$ws.Cells.Item(7, 6).Value = 41
$ws.Cells.Item(14, 6).Value = 156
$ws.Cells.Item(22, 6).Value = 18
$ws.Cells.Item(24, 6).Value = 918
$ws.Cells.Item(26, 6).Value = 1026
$ws.Cells.Item(27, 6).Value = 155
$ws.Cells.Item(31, 6).Value = 140

$ws = $wb.Worksheets.Item("本地生活")
$ws.Cells.Item(6, 6).Value = 728
$ws.Cells.Item(7, 6).Value = 259
$ws.Cells.Item(8, 6).Value = 127
$ws.Cells.Item(9, 6).Value = 1953
$ws.Cells.Item(10, 6).Value = 2960

$ws = $wb.Worksheets.Item("全部类型")
$ws.Cells.Item(5, 6).Value = 728
$ws.Cells.Item(7, 6).Value = 8723
$ws.Cells.Item(8, 6).Value = 259
$ws.Cells.Item(9, 6).Value = 127
$ws.Cells.Item(11, 6).Value = 1953
$ws.Cells.Item(12, 6).Value = 2960
$ws.Cells.Item(16, 6).Value = 96
$ws.Cells.Item(17, 6).Value = 1240
$ws.Cells.Item(18, 6).Value = 156
$ws.Cells.Item(20, 6).Value = 591
$ws.Cells.Item(21, 6).Value = 104
$ws.Cells.Item(22, 6).Value = 278
$ws.Cells.Item(23, 6).Value = 1402
$ws.Cells.Item(24, 6).Value = 1305
$ws.Cells.Item(25, 6).Value = 1324
$ws.Cells.Item(26, 6).Value = 205
$ws.Cells.Item(27, 6).Value = 57
$ws.Cells.Item(28, 6).Value = 273
$ws.Cells.Item(30, 6).Value = 5
$ws.Cells.Item(33, 6).Value = 918
$ws.Cells.Item(34, 6).Value = 199
$ws.Cells.Item(36, 6).Value = 174
$ws.Cells.Item(37, 6).Value = 156
$ws.Cells.Item(39, 6).Value = 596
$ws.Cells.Item(41, 6).Value = 664
$ws.Cells.Item(44, 6).Value = 195
$ws.Cells.Item(48, 6).Value = 38
